# Insert a new weekly price record for Arándano (blue) / Vega Modelo de Temuco
# as row 47, shifting all subsequent rows (old 47-132) down to (48-133).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 47 (pushes existing rows 47..132 down to 48..133).
$ws.Rows("47:47").Insert()

# Populate the newly inserted row 47 with the new data point.
$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value = "La Araucanía"
$ws.Cells.Item(47, 4).Value = 44935
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(47, 6).Value = "Fruta"
$ws.Cells.Item(47, 7).Value = 100101
$ws.Cells.Item(47, 8).Value = "Berries"
$ws.Cells.Item(47, 9).Value = 100101001
$ws.Cells.Item(47, 10).Value = "Arándano (blue)"
$ws.Cells.Item(47, 11).Value = "Sin especificar"
$ws.Cells.Item(47, 12).Value = "Primera"
$ws.Cells.Item(47, 13).Value = 500
$ws.Cells.Item(47, 14).Value = 1300
$ws.Cells.Item(47, 15).Value = 1300
$ws.Cells.Item(47, 16).Value = 1300
$ws.Cells.Item(47, 17).Value = "`$/kilo"
$ws.Cells.Item(47, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(47, 19).Value = 1300
$ws.Cells.Item(47, 20).Value = 1
